$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lookup name cell
$ws.Range("G2").Value = "BRIAN"

# Add helper boolean formulas in columns F:H
$ws.Range("F8").Formula = '=OR(ISBLANK($B2),ISBLANK($C2))'
$ws.Range("G8:H8").Formula = '=OR(ISBLANK($B2),ISBLANK($C2))'
$ws.Range("F9:H17").Formula = '=OR(ISBLANK($B3),ISBLANK($C3))'

# Conditional formatting on A2:C11
$rng = $ws.Range("A2:C11")

# -- history of dxf creation or edits that left orphaned dxf entries --
$fc = $rng.FormatConditions.Add(2, 0, '=OR(ISBLANK($B2),ISBLANK($C2))')
$fc.Interior.Color = 65535
$null = $rng.FormatConditions.Delete()

$fc = $rng.FormatConditions.Add(2, 0, '=$A2=$G$2')
$fc.Interior.Color = 15773696
$null = $rng.FormatConditions.Delete()

# -- kept rules --
$fcYellow = $rng.FormatConditions.Add(2, 0, '=OR(ISBLANK($B2),ISBLANK($C2))')
$fcYellow.Interior.Color = 65535

$fcBlue = $rng.FormatConditions.Add(2, 0, '=$A2=$G$2')
$fcBlue.Interior.Color = 15773696

# -- more orphaned dxf entries --
$fc = $rng.FormatConditions.Add(2, 0, '=$A2=$G$2')
$fc.Interior.Color = 15773696
$null = $fc.Delete()

$fc = $rng.FormatConditions.Add(2, 0, '=$A2=$G$2')
$fc.Interior.Color = 15773696
$null = $fc.Delete()

$fc = $rng.FormatConditions.Add(2, 0, '=$A2=$G$2')
$fc.Interior.Color = 15773696
$null = $fc.Delete()

# Selection
$null = $ws.Range("B4").Select()
